$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.809.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.471.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("E9").Value = "  +5.92%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "68.717.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000171"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "338.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0₃0828"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "430.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "130.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0721"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.488"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.565"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0919"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("E49").Value = "  -6.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("E51").Value = "  -5.56%  "
